$wb = $excel.ActiveWorkbook
$wsAll = $wb.Worksheets.Item("All")

# --- Header row updates ---
$wsAll.Range("D1").Value = "Organism"
$wsAll.Range("F1").Value = $wsAll.Range("G1").Value()
$wsAll.Range("G1").ClearContents()

# --- Uppercase gene names in column A (rows 2-43) ---
for ($r = 2; $r -le 43; $r++) {
    $cell = $wsAll.Cells.Item($r, 1)
    $cell.Value = $cell.Value().ToUpper()
}

# --- Move "Opmerkingen" (remarks) from column G into column F, clear G ---
$wsAll.Range("F10").Value = $wsAll.Range("G10").Value()
$wsAll.Range("G10").ClearContents()
$wsAll.Range("F14").Value = $wsAll.Range("G14").Value()
$wsAll.Range("G14").ClearContents()

# --- Replace remaining "Essentiality according to SCMD2" values in column F ---
# Rows whose SGD/PB essentiality (E) disagreed with SCMD2 (F): adopt SCMD2's
# "Essential" verdict into E, and mark F as unresolved ("?")
$wsAll.Range("E18").Value = "Essential"
$wsAll.Range("F18").Value = "?"
$wsAll.Range("F29").Value = "?"
$wsAll.Range("E33").Value = "Essential"
$wsAll.Range("F33").Value = "?"

# Rows whose SCMD2 essentiality column is simply cleared
$wsAll.Range("F2").ClearContents()
$wsAll.Range("F3").ClearContents()
$wsAll.Range("F4").ClearContents()
$wsAll.Range("F5").ClearContents()
$wsAll.Range("F6").ClearContents()
$wsAll.Range("F7").ClearContents()
$wsAll.Range("F8").ClearContents()
$wsAll.Range("F9").ClearContents()
$wsAll.Range("F11").ClearContents()
$wsAll.Range("F12").ClearContents()
$wsAll.Range("F13").ClearContents()
$wsAll.Range("F16").ClearContents()
$wsAll.Range("F17").ClearContents()
$wsAll.Range("F19").ClearContents()
$wsAll.Range("F20").ClearContents()
$wsAll.Range("F21").ClearContents()
$wsAll.Range("F22").ClearContents()
$wsAll.Range("F23").ClearContents()
$wsAll.Range("F24").ClearContents()
$wsAll.Range("F25").ClearContents()
$wsAll.Range("F26").ClearContents()
$wsAll.Range("F27").ClearContents()
$wsAll.Range("F28").ClearContents()
$wsAll.Range("F30").ClearContents()
$wsAll.Range("F32").ClearContents()
$wsAll.Range("F34").ClearContents()
$wsAll.Range("F35").ClearContents()
$wsAll.Range("F36").ClearContents()
$wsAll.Range("F37").ClearContents()
$wsAll.Range("F38").ClearContents()
$wsAll.Range("F39").ClearContents()
$wsAll.Range("F42").ClearContents()
$wsAll.Range("F43").ClearContents()

# --- Sheet "Absent": uppercase species/gene identifiers in column A ---
$wsAbsent = $wb.Worksheets.Item("Absent")
for ($r = 1; $r -le 6; $r++) {
    $cell = $wsAbsent.Cells.Item($r, 1)
    $cell.Value = $cell.Value().ToUpper()
}
